$wb = $excel.ActiveWorkbook

# Sheet: 展览
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 2277
$ws1.Range("F3").Value = 348
$ws1.Range("F4").Value = 182
$ws1.Range("F5").Value = 186
$ws1.Range("F6").Value = 347
$ws1.Range("F8").Value = 711
$ws1.Range("F9").Value = 518
$ws1.Range("F10").Value = 675
$ws1.Range("F11").Value = 371
$ws1.Range("F13").Value = 362
$ws1.Range("F14").Value = 974
$ws1.Range("F15").Value = 5642
$ws1.Range("F16").Value = 182
$ws1.Range("F17").Value = 16
$ws1.Range("F18").Value = 33
$ws1.Range("F19").Value = 252
$ws1.Range("F21").Value = 112
$ws1.Range("F23").Value = 94
$ws1.Range("F25").Value = 263
$ws1.Range("F26").Value = 107

# Sheet: 演出
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F6").Value = 183
$ws2.Range("F8").Value = 2801
$ws2.Range("F13").Value = 26
$ws2.Range("F14").Value = 108
$ws2.Range("F16").Value = 2529

# Sheet: 本地生活
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F3").Value = 45
$ws3.Range("F4").Value = 406

# Sheet: 全部类型
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 45
$ws4.Range("F6").Value = 2277
$ws4.Range("F7").Value = 406
$ws4.Range("F8").Value = 348
$ws4.Range("F9").Value = 182
$ws4.Range("F10").Value = 186
$ws4.Range("F11").Value = 347
$ws4.Range("F15").Value = 183
$ws4.Range("F17").Value = 711
$ws4.Range("F18").Value = 518
$ws4.Range("F19").Value = 675
$ws4.Range("F20").Value = 371
$ws4.Range("F22").Value = 362
$ws4.Range("F23").Value = 974
$ws4.Range("F24").Value = 5646
$ws4.Range("F26").Value = 2801
$ws4.Range("F30").Value = 182
$ws4.Range("F31").Value = 16
$ws4.Range("F32").Value = 33
$ws4.Range("F34").Value = 26
$ws4.Range("F35").Value = 252
$ws4.Range("F37").Value = 112
$ws4.Range("F39").Value = 108
$ws4.Range("F41").Value = 94
$ws4.Range("F43").Value = 263
$ws4.Range("F44").Value = 107
$ws4.Range("F45").Value = 2529
